$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.9215594866898087, 1.7331519473645112, -1.5469261442825584, 2.5449011747400863)"
$ws.Range("C2").Value = "NIG(2.0703748312614274, 1.6968808133207323, 3.230016032244919, 6.643863361639738)"
$ws.Range("D2").Value = "NIG(0.8395308541281241, 0.5467952768774937, 1.3115356125259146, 3.11366034188084)"
$ws.Range("E2").Value = "NIG(1.174633486569486, 0.8909627719899003, 3.394944577059021, 6.109818489204002)"
